$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Germany")
$ws1.Activate()
$ws1.Range("A3").Select()

$ws2 = $wb.Worksheets.Item("Belgium")
$ws2.Activate()
$ws2.Range("A8:A35").Select()

Write-Output "done"
